$d = $word.ActiveDocument

# --- Locate the "BIBtex:" paragraph -------------------------------------
$bibParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "BIBtex:`r") {
        $bibParaIndex = $i
        break
    }
}
if ($bibParaIndex -eq -1) {
    throw "Could not find the 'BIBtex:' paragraph"
}

$bibPara = $d.Paragraphs.Item($bibParaIndex)
$bibRange = $bibPara.Range

# "BIBtex:" -> "BIBtex: " (add the trailing space) inside the existing run.
$textRange = $d.Range($bibRange.Start, $bibRange.End - 1)
$textRange.Text = "BIBtex: "

# Append the fenced-code opening marker right after it.
$bibPara2 = $d.Paragraphs.Item($bibParaIndex)
$bibRange2 = $bibPara2.Range
$insertionPoint = $d.Range($bibRange2.End - 1, $bibRange2.End - 1)
$insertionPoint.InsertAfter("``````")

# The newly inserted text has identical run formatting to "BIBtex: " so it
# gets folded back into the same run; force it into its own run (matching
# the target markup, two sibling <w:r> elements) by toggling Bold off/on
# across just that span.
$bibPara3 = $d.Paragraphs.Item($bibParaIndex)
$bibRange3 = $bibPara3.Range
$newRunRange = $d.Range($bibRange3.End - 4, $bibRange3.End - 1)
$newRunRange.Font.Bold = $false
$newRunRange.Font.Bold = $true

# --- Locate the blank "Quotations" paragraph immediately after the
#     closing "}" of the BibTeX entry, and add the closing fence there --
$closeBraceIndex = -1
for ($i = $bibParaIndex; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "}`r") {
        $closeBraceIndex = $i
        break
    }
}
if ($closeBraceIndex -eq -1) {
    throw "Could not find the closing '}' paragraph of the BibTeX entry"
}

$fenceParaIndex = $closeBraceIndex + 1
$fencePara = $d.Paragraphs.Item($fenceParaIndex)
if ($fencePara.Range.Text -ne "`r") {
    throw "Expected an empty paragraph right after the closing '}'"
}

# Pin the (until now style-inherited) right indent / hanging indent
# explicitly onto this paragraph -> <w:ind w:right="567" w:hanging="0"/>.
$fencePara.Range.ParagraphFormat.RightIndent = 28.35
$fencePara.Range.ParagraphFormat.FirstLineIndent = -0.001

# Insert the closing fence text into the (until now empty) run.
$fencePara2 = $d.Paragraphs.Item($fenceParaIndex)
$fenceRange2 = $fencePara2.Range
$fenceInsertionPoint = $d.Range($fenceRange2.Start, $fenceRange2.Start)
$fenceInsertionPoint.InsertAfter("``````")

Write-Output "Applied BIBtex fenced-code-block edits."
